$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.621.42"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "'1.596.54"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'211.17"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.247"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'19.51"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'0.0839"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "'1.820.22"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'1.576.66"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'26.601.39"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'0.0₃0739"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'208.60"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "'6.99"
$ws.Range("E21").Value = "  +4.85%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "'145.40"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'15.27"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "'0.0511"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").Value = "'1.280.83"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").Value = "'0.622"
$ws.Range("E35").Value = "  -6.51%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'0.839"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  +18.67%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").Value = "'64.19"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").Value = "'1.733.03"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'90.05"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "'7.45"
$ws.Range("E51").Value = "  -1.07%  "
